$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new header values
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the formatting from an existing header cell (bold, bordered, centered)
# onto the new header cells so they reuse the same style.
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New boolean columns (all FALSE) for rows 2-4
$ws.Range("F2:H4").Value = $false
